# Update database (shift quarterly columns by one, adding newest quarter,
# dropping oldest) and adjust column widths / row heights to match the
# newly-saved workbook's layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 8): quarter labels, shifted left by one column with
#    a new quarter appended in column M.
# ---------------------------------------------------------------------
$periodLabels = @(
    "فصل چهارم منتهی به 1399/08",
    "فصل اول منتهی به 1399/11",
    "فصل دوم منتهی به 1400/02",
    "فصل سوم منتهی به 1400/05",
    "فصل چهارم منتهی به 1400/08",
    "فصل اول منتهی به 1400/11",
    "فصل دوم منتهی به 1401/02",
    "فصل سوم منتهی به 1401/05",
    "فصل چهارم منتهی به 1401/08",
    "فصل اول منتهی به 1401/11"
)
for ($i = 0; $i -lt $periodLabels.Length; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value = $periodLabels[$i]
}

# ---------------------------------------------------------------------
# 2. Publish-date row (row 9), shifted left by one column with a new
#    date appended in column M.
# ---------------------------------------------------------------------
$publishDates = @(
    "1400-12-09 (9)",
    "1400-12-28 (2)",
    "1401-05-12 (4)",
    "1401-06-30 (2)",
    "1401-12-13 (9)",
    "1401-12-28 (2)",
    "1401-05-12 (2)",
    "1401-06-30",
    "1401-12-28 (3)",
    "1401-12-28"
)
for ($i = 0; $i -lt $publishDates.Length; $i++) {
    $ws.Cells.Item(9, 4 + $i).Value = $publishDates[$i]
}

# ---------------------------------------------------------------------
# 3. Financial data rows (11-27), each shifted left by one column and a
#    new value (or "-") appended in column M. Rows that are entirely
#    "-" or entirely 0 are unaffected by the shift and left untouched.
# ---------------------------------------------------------------------
$rowData = @{
    11 = @(8536, 9228, 12272, 8973, 10836, 10313, 11348, 12185, 13428, 8815)
    12 = @(-6052, -5647, -7449, -5969, -10429, -6902, -7691, -9072, -11229, -6345)
    13 = @(2483, 3581, 4823, 3004, 407, 3411, 3657, 3113, 2199, 2470)
    14 = @(-264, -115, -564, 197, -591, -148, -736, -51, -191, -235)
    16 = @(182, 9, -52, 41, -109, "-", 274, 62, -170, 356)
    17 = @(2402, 3475, 4207, 3242, -293, 3264, 3195, 3125, 1838, 2591)
    18 = @(-530, -588, -478, -184, -762, -503, -496, -629, -542, -475)
    19 = @(152, 1, 154, -84, 1586, "-", 142, 84, -150, "-")
    20 = @(2024, 2888, 3883, 2975, 530, 2760, 2842, 2580, 1146, 2117)
    21 = @(-280, -561, -998, -355, 1673, "-", "-", "-", "-", "-")
    22 = @(1745, 2327, 2886, 2620, 2203, 2760, 2842, 2580, 1146, 2117)
    24 = @(1745, 2327, 2886, 2620, 2203, 2760, 2842, 2580, 1146, 2117)
    26 = @(35283, 37693, 39614, 37806, 33828, 32654, 34008, 29827, 28787, 22839)
}

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, 4 + $i).Value = $values[$i]
    }
}

# ---------------------------------------------------------------------
# 4. Column widths - re-balanced by one column, matching the resave.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 30.166666666666668    # D -> 31
$ws.Columns.Item(5).ColumnWidth = 28.166666666666668    # E -> 29
$ws.Columns.Item(6).ColumnWidth = 28.166666666666668    # F -> 29
$ws.Columns.Item(7).ColumnWidth = 28.166666666666668    # G -> 29
$ws.Columns.Item(8).ColumnWidth = 30.166666666666668    # H -> 31
$ws.Columns.Item(9).ColumnWidth = 28.166666666666668    # I -> 29
$ws.Columns.Item(10).ColumnWidth = 28.166666666666668   # J -> 29
$ws.Columns.Item(11).ColumnWidth = 28.166666666666668   # K -> 29
$ws.Columns.Item(12).ColumnWidth = 30.166666666666668   # L -> 31
$ws.Columns.Item(13).ColumnWidth = 28.166666666666668   # M -> 29

# ---------------------------------------------------------------------
# 5. Row heights - default row height / title rows changed slightly on
#    resave (Calibri metrics recalculation).
# ---------------------------------------------------------------------
$ws.StandardHeight = 14.4
$ws.Rows.Item(2).RowHeight = 15.6
$ws.Rows.Item(5).RowHeight = 40.8
$ws.Rows.Item(6).RowHeight = 40.8
$ws.Rows.Item(8).RowHeight = 15.6

$wb.Save()
